# Add a "Std Dev" column next to the existing "Avg" column.
#
# F2 becomes a header cell (text "Std Dev"), formatted like the other
# header cells in row 2 (B2/C2/E2).
# F3 becomes a formula cell computing STDEV(B3:B5), formatted like the
# other formula/value cells in row 3 (C3/E3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: header label "Std Dev" -------------------------------------
# Copy the formatting of the neighboring header cell E2 ("Avg") onto
# F2, then set its text value.
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F2").Value = "Std Dev"

# --- F3: formula =STDEV(B3:B5) ---------------------------------------
# Copy the formatting of the neighboring formula cell E3 (AVERAGE)
# onto F3, then set its formula.
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F3").Formula = "=STDEV(B3:B5)"

$excel.CutCopyMode = 0
